$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "cluster_length" column header in G1, matching F1's header formatting
$ws.Range("G1").Value = "cluster_length"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill G2 first (non-shared anchor formula), then G3:G66 and G67:G83 as separate fill groups
$ws.Range("G2").Formula = "=COUNTIF(C:C, C2)"
$ws.Range("G3:G66").Formula = "=COUNTIF(C:C, C3)"
$ws.Range("G67:G83").Formula = "=COUNTIF(C:C, C67)"

# Match the original selection state (H1 instead of H3)
$ws.Range("H1").Select()
